$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 30,6
$arr[0,0] = 'Manager in Training (MIT) Program - Middle East'
$arr[0,1] = 'Louis Vuitton'
$arr[0,2] = 'Dubai'
$arr[0,3] = 'As part of Louis Vuitton’s Manager in Training (MIT) Program, you will embark upon a dynamic and in-depth journey of discovery through 4 phases: Client advisor…'
$arr[0,4] = 'Posted4 days ago'
$arr[0,5] = 'https://ae.indeed.com/rc/clk?jk=7141066c9fc08d8a&fccid=1807e5727f702882&vjs=3'
$arr[1,0] = 'General Manager - Luxury Resort'
$arr[1,1] = 'Michael Page AE'
$arr[1,2] = 'Dubai'
$arr[1,3] = ' At least 10 years of experience in the executive management of operations, sales and marketing, food and beverage or related professional area of luxury…'
$arr[1,4] = 'Posted30+ days ago'
$arr[1,5] = 'https://ae.indeed.com/rc/clk?jk=f15a7ddaf4d3bef8&fccid=77087bd1709a8148&vjs=3'
$arr[2,0] = 'Assistant Department Manager LG'
$arr[2,1] = 'Christian Dior Couture'
$arr[2,2] = 'Dubai'
$arr[2,3] = 'Follow up with the Department managers to ensure optimization of the stock level usage.'
$arr[2,4] = 'Posted25 days ago'
$arr[2,5] = 'https://ae.indeed.com/rc/clk?jk=5b0e5213019a6fd8&fccid=6a2be34af774e2bb&vjs=3'
$arr[3,0] = 'Duty Manager'
$arr[3,1] = 'Anantara'
$arr[3,2] = 'Abu Dhabi'
$arr[3,3] = 'Should be a strong and experienced hotel ambassador who excels at maximizing the experience and satisfaction of the hotel guests.'
$arr[3,4] = 'Posted16 days ago'
$arr[3,5] = 'https://ae.indeed.com/rc/clk?jk=850fad99487d4240&fccid=2af5b593acca1f1f&vjs=3'
$arr[4,0] = 'Assistant Manager - Dubai Mall'
$arr[4,1] = 'Tiffany & Co'
$arr[4,2] = 'Dubai'
$arr[4,3] = 'A minimum three year’s retail leadership experience within a luxury retail environment.'
$arr[4,4] = 'Posted30+ days ago'
$arr[4,5] = 'https://ae.indeed.com/rc/clk?jk=bfca47ad77a2ebea&fccid=329e52936d287237&vjs=3'
$arr[5,0] = 'Regional Retail Manager - Premium Luxury Brands'
$arr[5,1] = 'Chalhoub Group'
$arr[5,2] = 'Dubai'
$arr[5,3] = 'Working as a Retail Manager you will be responsible for ensuring that Store Directors/Managers maximise the commercial performance of their retail shops and…'
$arr[5,4] = 'Posted30+ days ago'
$arr[5,5] = 'https://ae.indeed.com/rc/clk?jk=3d46626a0e65ba12&fccid=01f47b3f00b281a4&vjs=3'
$arr[6,0] = 'Boutique Manager - Al Ain Mall'
$arr[6,1] = 'Chanel'
$arr[6,2] = 'Abu Dhabi'
$arr[6,3] = 'You have at least 5 years’ experience in boutique management in a luxury retail environment preferably in fragrance and beauty.'
$arr[6,4] = 'Posted12 days ago'
$arr[6,5] = 'https://ae.indeed.com/rc/clk?jk=15c6e1d5546301f2&fccid=20f48cf7726b0510&vjs=3'
$arr[7,0] = 'Showroom Manager – Africa'
$arr[7,1] = 'Danube Group'
$arr[7,2] = 'Dubai'
$arr[7,3] = 'This luxury furniture is made with the highest quality, all-natural, material.'
$arr[7,4] = 'Posted13 days ago'
$arr[7,5] = 'https://ae.indeed.com/rc/clk?jk=f24791aa7942acd7&fccid=a6bec4b75595280e&vjs=3'
$arr[8,0] = 'Retail Design Project Manager'
$arr[8,1] = 'L''Oreal'
$arr[8,2] = 'Remote in Dubai'
$arr[8,3] = 'Within the Retail design and visual merchandising department, the Retail design project manager will be in charge of the implementation of the permanent POS …'
$arr[8,4] = 'Posted9 days ago'
$arr[8,5] = 'https://ae.indeed.com/rc/clk?jk=1d7504e169401aed&fccid=ec8a0bd60be9a6f5&vjs=3'
$arr[9,0] = 'Retail Manager - Guerlain'
$arr[9,1] = 'Chalhoub Group'
$arr[9,2] = 'Dubai'
$arr[9,3] = 'The Regional Retail Manager is responsible for sales objective achievement and outstanding omnichannel client experience across the retail network (boutiques…'
$arr[9,4] = 'Posted3 days ago'
$arr[9,5] = 'https://ae.indeed.com/rc/clk?jk=c183b9675d782014&fccid=01f47b3f00b281a4&vjs=3'
$arr[10,0] = 'Market Development Manager - Dubai, UAE'
$arr[10,1] = 'Argyll Scott MY'
$arr[10,2] = 'Dubai'
$arr[10,3] = 'Up to UAE Dirhams384000 per annum per annum.'
$arr[10,4] = 'Posted5 days ago'
$arr[10,5] = 'https://ae.indeed.com/rc/clk?jk=cc0a01bd82e90fc5&fccid=91277dc7f9bcbc55&vjs=3'
$arr[11,0] = 'Senior Regional Manager - Luxury Fragrances - MEA Coverage'
$arr[11,1] = 'Michael Page AE'
$arr[11,2] = 'Dubai'
$arr[11,3] = 'Reporting to the Managing Director, this person will be responsible for;'
$arr[11,4] = 'Posted30+ days ago'
$arr[11,5] = 'https://ae.indeed.com/rc/clk?jk=f14fd6c319542ff6&fccid=77087bd1709a8148&vjs=3'
$arr[12,0] = 'Team Manager - Client Service Center'
$arr[12,1] = 'Louis Vuitton'
$arr[12,2] = 'Dubai'
$arr[12,3] = 'Manage and motivate the team to drive business: create a positive and harmonious work environment, foster cooperation within the team and between managers.'
$arr[12,4] = 'Posted20 days ago'
$arr[12,5] = 'https://ae.indeed.com/rc/clk?jk=add0a2677658044d&fccid=1807e5727f702882&vjs=3'
$arr[13,0] = 'Night Manager'
$arr[13,1] = 'Marriott International, Inc'
$arr[13,2] = 'Dubai'
$arr[13,3] = 'Assists operations manager in processing employee payroll weekly.'
$arr[13,4] = 'Posted13 days ago'
$arr[13,5] = 'https://ae.indeed.com/rc/clk?jk=ffbdce9fa2d04545&fccid=0b6c496064ecd79a&vjs=3'
$arr[14,0] = 'Assistant EBC Manager'
$arr[14,1] = 'Marriott International, Inc'
$arr[14,2] = 'Dubai'
$arr[14,3] = 'Contact appropriate individual or department (e.g., Sales, Data Administration, Accounting) as necessary to resolve guest calls, requests, or problems.'
$arr[14,4] = 'Posted2 days ago'
$arr[14,5] = 'https://ae.indeed.com/rc/clk?jk=1b0e4e550137f8a7&fccid=0b6c496064ecd79a&vjs=3'
$arr[15,0] = 'MEA Logistics Manager'
$arr[15,1] = 'Luxury Goods International (L.G.I) SA (Branch)'
$arr[15,2] = 'Dubai'
$arr[15,3] = 'Based in Dubai, the MEA Logistics Manager supports the development of the Kering Brands’ business by implementing and running best-in-class logistics solutions…'
$arr[15,4] = 'Posted30+ days ago'
$arr[15,5] = 'https://ae.indeed.com/rc/clk?jk=3de3673c2ca392d7&fccid=dd616958bd9ddc12&vjs=3'
$arr[16,0] = 'Sales Delgate/ Junior Area Retail Supervisor'
$arr[16,1] = 'SSC Perfumes & Cosmetics'
$arr[16,2] = 'Dubai'
$arr[16,3] = 'LVMH P&C Middle East oversees an area of 47 countries in the Middle East, Europe, Africa, and India with subsidiaries and agents’ structures.'
$arr[16,4] = 'Posted25 days ago'
$arr[16,5] = 'https://ae.indeed.com/rc/clk?jk=57d8e49629b99c35&fccid=5ad360b814db19a5&vjs=3'
$arr[17,0] = 'Senior Project Manager - Store Renovation - Level Shoes'
$arr[17,1] = 'Chalhoub Group'
$arr[17,2] = 'Remote in Dubai'
$arr[17,3] = 'Excellent time management with capabilities to multitask with a strong understanding of core manager duties.'
$arr[17,4] = 'Posted18 days ago'
$arr[17,5] = 'https://ae.indeed.com/rc/clk?jk=1e399b9031998ae9&fccid=01f47b3f00b281a4&vjs=3'
$arr[18,0] = 'Clienteling Manager – CELINE Dubai Mall'
$arr[18,1] = 'Chalhoub Group'
$arr[18,2] = 'Dubai'
$arr[18,3] = 'Support store managers in managing their teams on customer issues.'
$arr[18,4] = 'Posted20 days ago'
$arr[18,5] = 'https://ae.indeed.com/rc/clk?jk=72cd01f90905f251&fccid=01f47b3f00b281a4&vjs=3'
$arr[19,0] = 'Assistant Manager - Dubai Mall'
$arr[19,1] = 'Tiffany & Co.'
$arr[19,2] = 'Dubai'
$arr[19,3] = 'A minimum three year’s retail leadership experience within a luxury retail environment.'
$arr[19,4] = 'Posted30+ days ago'
$arr[19,5] = 'https://ae.indeed.com/rc/clk?jk=83b998ce18b6cb88&fccid=329e52936d287237&vjs=3'
$arr[20,0] = 'Private Client Manager'
$arr[20,1] = 'FARFETCH'
$arr[20,2] = 'Dubai'
$arr[20,3] = 'Through a variety of engagement and selling activities, you''ll ensure a seamless, exceptional luxury shopping experience for our Private Client customers.'
$arr[20,4] = 'Posted30+ days ago'
$arr[20,5] = 'https://ae.indeed.com/rc/clk?jk=b098d592cdae6db1&fccid=9f5fb2f8ae2fcc49&vjs=3'
$arr[21,0] = 'Private Client Assistant Stylist'
$arr[21,1] = 'FARFETCH'
$arr[21,2] = 'Dubai'
$arr[21,3] = 'We''re looking for someone with strong customer service skills and knowledge of the luxury fashion market who is driven to achieve targets.'
$arr[21,4] = 'Posted2 days ago'
$arr[21,5] = 'https://ae.indeed.com/rc/clk?jk=0123577ae0ac18b3&fccid=9f5fb2f8ae2fcc49&vjs=3'
$arr[22,0] = 'Area Retail Manager – Luxury Brand'
$arr[22,1] = 'Charterhouse Consultancy PTE Ltd'
$arr[22,2] = 'Dubai'
$arr[22,3] = 'The main purpose of this role is to meet the business objectives, drive the performance and maximize the profitability of all the assigned stores, recruit,…'
$arr[22,4] = 'Posted30+ days ago'
$arr[22,5] = 'https://ae.indeed.com/rc/clk?jk=058b3df0807cd28e&fccid=4749f34298e2e970&vjs=3'
$arr[23,0] = 'E-Commerce Manager'
$arr[23,1] = 'coty'
$arr[23,2] = 'Dubai'
$arr[23,3] = 'You will be responsible for developing and overseeing the brands online sales, seamless user experience and customer journey on the web.'
$arr[23,4] = 'Posted30+ days ago'
$arr[23,5] = 'https://ae.indeed.com/rc/clk?jk=f2ede358b9f91676&fccid=7ad46606e93080a6&vjs=3'
$arr[24,0] = 'Quality Assurance Manager - Emirates Academy of Hospitality Management'
$arr[24,1] = 'Jumeirah'
$arr[24,2] = 'Dubai'
$arr[24,3] = 'The Quality Assurance Manager is directly responsible to the Dean, EAHM, for the accurate and effective reporting of institutional data, academic performance.'
$arr[24,4] = 'Posted30+ days ago'
$arr[24,5] = 'https://ae.indeed.com/rc/clk?jk=17a460bfdfebc7fd&fccid=8015c178481add76&vjs=3'
$arr[25,0] = 'Junior Area Retail Manager - Guerlain'
$arr[25,1] = 'Chalhoub Group'
$arr[25,2] = 'Dubai'
$arr[25,3] = 'As Junior Area Retail Manager, you will be responsible for supporting the Area Retail Manager within Guerlain for Travel Retail for the region.'
$arr[25,4] = 'Posted20 days ago'
$arr[25,5] = 'https://ae.indeed.com/rc/clk?jk=8793ce2c66d0caff&fccid=01f47b3f00b281a4&vjs=3'
$arr[26,0] = 'Recreation Manager'
$arr[26,1] = 'Anantara'
$arr[26,2] = 'Dubai'
$arr[26,3] = 'Oversee the smooth running of the Sports & Recreations department on a day to day basis to achieve the highest possible levels of guest satisfaction.'
$arr[26,4] = 'Posted27 days ago'
$arr[26,5] = 'https://ae.indeed.com/rc/clk?jk=41d202b8068a372c&fccid=2af5b593acca1f1f&vjs=3'
$arr[27,0] = 'Strategy Manager'
$arr[27,1] = 'Chalhoub Group'
$arr[27,2] = 'Dubai'
$arr[27,3] = 'The Strategy Office provides the organisation with an overall compass, based on the Group’s vision and directions previously defined.'
$arr[27,4] = 'Posted30+ days ago'
$arr[27,5] = 'https://ae.indeed.com/rc/clk?jk=3794f02a2013b247&fccid=01f47b3f00b281a4&vjs=3'
$arr[28,0] = 'Senior Visual Merchandiser - Luxury Fashion - Abu Dhabi'
$arr[28,1] = 'Chalhoub Group'
$arr[28,2] = 'Abu Dhabi'
$arr[28,3] = 'Communicate with retail manager on impact of new merchandising on sales.'
$arr[28,4] = 'Posted6 days ago'
$arr[28,5] = 'https://ae.indeed.com/rc/clk?jk=df472563e2b34dd6&fccid=01f47b3f00b281a4&vjs=3'
$arr[29,0] = 'Property Manager'
$arr[29,1] = 'Asayel Investment'
$arr[29,2] = 'Abu Dhabi'
$arr[29,3] = 'Budgeting. • Real estate financial analysis.'
$arr[29,4] = 'Posted30+ days ago'
$arr[29,5] = 'https://ae.indeed.com/rc/clk?jk=16b9779594f116f3&fccid=f697fec203d7893a&vjs=3'

$ws.Range("A2:F31").Value = $arr